$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VLAN_POOL")

# Replace the numeric VLAN range values with text-based vlan pool names
$ws.Range("C2").Value = "vlan-1024"
$ws.Range("D2").Value = "vlan-1034"

# Activate the sheet and move the selection to D2, matching the saved cursor position
$ws.Activate()
$ws.Range("D2").Select()
